$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the crypto price-refresh diff.
# Cells whose new value parses as a plain number need NumberFormat "@"
# (Text) applied first so Excel keeps storing them as literal strings
# (matching the original inlineStr cell contents, e.g. "0.999", "0.190").

$ws.Range("D2").Value = '42.576.92'
$ws.Range("E2").Value = '  -1.78%  '
$ws.Range("D3").Value = '2.516.05'
$ws.Range("E3").Value = '  -3.42%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '308.09'
$ws.Range("E5").Value = '  -2.79%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '100.62'
$ws.Range("E6").Value = '  +2.19%  '
$ws.Range("E7").Value = '  -1.74%  '
$ws.Range("E8").Value = '  +0.12%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.525'
$ws.Range("E9").Value = '  -3.10%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.94'
$ws.Range("E10").Value = '  -0.50%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0801'
$ws.Range("E11").Value = '  -2.07%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.29'
$ws.Range("E12").Value = '  -3.78%  '
$ws.Range("E13").Value = '  -0.21%  '
$ws.Range("D14").Value = '2.907.82'
$ws.Range("E14").Value = '  -3.26%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.58'
$ws.Range("E15").Value = '  +2.13%  '
$ws.Range("D16").Value = '2.517.99'
$ws.Range("E16").Value = '  -4.78%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.806'
$ws.Range("E17").Value = '  -5.36%  '
$ws.Range("D18").Value = '42.545.66'
$ws.Range("E18").Value = '  -2.06%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.71'
$ws.Range("E19").Value = '  -2.44%  '
$ws.Range("D20").Value = '0.0₃0947'
$ws.Range("E20").Value = '  -2.59%  '
$ws.Range("E21").Value = '  -5.50%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '69.44'
$ws.Range("E22").Value = '  -0.42%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '243.62'
$ws.Range("E23").Value = '  -4.50%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.88'
$ws.Range("E24").Value = '  -3.53%  '
$ws.Range("E25").Value = '  -2.89%  '
$ws.Range("E26").Value = '  +0.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.88'
$ws.Range("E27").Value = '  -5.40%  '
$ws.Range("E28").Value = '  -3.82%  '
$ws.Range("B29").Value = 'InjectiveProtocol'
$ws.Range("C29").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '38.86'
$ws.Range("E29").Value = '  -6.35%  '
$ws.Range("B30").Value = 'Cosmos'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '10.11'
$ws.Range("E30").Value = '  -2.21%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '155.81'
$ws.Range("E31").Value = '  -0.52%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.76'
$ws.Range("E32").Value = '  -2.41%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.77'
$ws.Range("E33").Value = '  +9.93%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0787'
$ws.Range("E34").Value = '  -3.24%  '
$ws.Range("E35").Value = '  -2.80%  '
$ws.Range("B36").Value = 'LidoDAOToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.18'
$ws.Range("E36").Value = '  -9.12%  '
$ws.Range("B37").Value = 'ARBITRUM'
$ws.Range("C37").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.02'
$ws.Range("E37").Value = '  -7.38%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.15'
$ws.Range("E38").Value = '  -3.69%  '
$ws.Range("E39").Value = '  -1.54%  '
$ws.Range("E40").Value = '  -0.36%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.26'
$ws.Range("E41").Value = '  +5.80%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '22.15'
$ws.Range("E42").Value = '  -3.38%  '
$ws.Range("E43").Value = '  +0.13%  '
$ws.Range("B44").Value = 'NEARProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.28'
$ws.Range("E44").Value = '  +0.33%  '
$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0298'
$ws.Range("E45").Value = '  -2.71%  '
$ws.Range("D46").Value = '1.981.47'
$ws.Range("E46").Value = '  -1.67%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.82'
$ws.Range("E47").Value = '  -2.14%  '
$ws.Range("D48").Value = '2.763.02'
$ws.Range("E48").Value = '  -3.15%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '80.09'
$ws.Range("E49").Value = '  -4.44%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.190'
$ws.Range("E50").Value = '  -3.08%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.39'
$ws.Range("E51").Value = '  -3.78%  '
